$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.0045489006823351
$ws.Range("C2").Value = 0.00985595147839272
$ws.Range("D2").Value = 0.0128885519332828
$ws.Range("E2").Value = 0.943896891584534
$ws.Range("F2").Value = 0.0151630022744503
$ws.Range("G2").Value = 0.980288097043215
$ws.Range("H2").Value = 0.967399545109932
$ws.Range("I2").Value = 0.0106141015921152
$ws.Range("J2").Value = 0.00606520090978014
$ws.Range("K2").Value = 0.00151630022744503
$ws.Range("L2").Value = 0.99090219863533
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0.000758150113722517
$ws.Range("R2").Value = 0.000758150113722517
$ws.Range("S2").Value = 0.0363912054586808
$ws.Range("T2").Value = 0.0295678544351782
$ws.Range("U2").Value = 0.00227445034116755
$ws.Range("V2").Value = 0.0151630022744503
$ws.Range("W2").Value = 0.0166793025018954
$ws.Range("X2").Value = 0.0045489006823351
$ws.Range("B3").Value = 0.0181956027293404
$ws.Range("C3").Value = 0.962850644427597
$ws.Range("D3").Value = 0.00682335102350265
$ws.Range("E3").Value = 0.00682335102350265
$ws.Range("F3").Value = 0.970432145564822
$ws.Range("G3").Value = 0.0174374526156179
$ws.Range("H3").Value = 0.00227445034116755
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0.00227445034116755
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0.998483699772555
$ws.Range("N3").Value = 0.00379075056861259
$ws.Range("O3").Value = 0.0045489006823351
$ws.Range("P3").Value = 0.00227445034116755
$ws.Range("Q3").Value = 0.99696739954511
$ws.Range("R3").Value = 0.989385898407885
$ws.Range("S3").Value = 0
$ws.Range("T3").Value = 0.954510993176649
$ws.Range("U3").Value = 0.020470053070508
$ws.Range("V3").Value = 0.00151630022744503
$ws.Range("W3").Value = 0
$ws.Range("X3").Value = 0
$ws.Range("B4").Value = 0.187263078089462
$ws.Range("C4").Value = 0.00303260045489007
$ws.Range("D4").Value = 0.0128885519332828
$ws.Range("E4").Value = 0.043972706595906
$ws.Range("F4").Value = 0.00379075056861259
$ws.Range("G4").Value = 0.00151630022744503
$ws.Range("H4").Value = 0.0250189537528431
$ws.Range("I4").Value = 0.988627748294162
$ws.Range("J4").Value = 0.991660348749052
$ws.Range("K4").Value = 0.995451099317665
$ws.Range("L4").Value = 0.00530705079605762
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0.000758150113722517
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 0.000758150113722517
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 0.000758150113722517
$ws.Range("S4").Value = 0.963608794541319
$ws.Range("T4").Value = 0.00379075056861259
$ws.Range("U4").Value = 0.000758150113722517
$ws.Range("V4").Value = 0.974981046247157
$ws.Range("W4").Value = 0.982562547384382
$ws.Range("X4").Value = 0.991660348749052
$ws.Range("B5").Value = 0.789992418498863
$ws.Range("C5").Value = 0.0242608036391205
$ws.Range("D5").Value = 0.967399545109932
$ws.Range("E5").Value = 0.00530705079605762
$ws.Range("F5").Value = 0.0106141015921152
$ws.Range("G5").Value = 0.000758150113722517
$ws.Range("H5").Value = 0.00530705079605762
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0.00303260045489007
$ws.Range("L5").Value = 0.00379075056861259
$ws.Range("M5").Value = 0.000758150113722517
$ws.Range("N5").Value = 0.995451099317665
$ws.Range("O5").Value = 0.995451099317665
$ws.Range("P5").Value = 0.99696739954511
$ws.Range("Q5").Value = 0.00227445034116755
$ws.Range("R5").Value = 0.0090978013646702
$ws.Range("S5").Value = 0
$ws.Range("T5").Value = 0.0121304018195603
$ws.Range("U5").Value = 0.976497346474602
$ws.Range("V5").Value = 0.00758150113722517
$ws.Range("W5").Value = 0
$ws.Range("X5").Value = 0.00303260045489007
